# Demo_Database.xlsx update
# - Insert a new "ID_User" attribute row into the User table
# - Rename the worksheets (drop the "Table_" prefix)
# - Switch the active/selected sheet from User to Comment, updating selections

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new row describing the "ID_User" attribute right below the
#    header row of the Table_User sheet, shifting the existing rows down.
# ---------------------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("Table_User")

$wsUser.Rows.Item(2).Insert()
$wsUser.Cells.Item(2, 1).Value = 1
$wsUser.Cells.Item(2, 2).Value = "ID_User"
$wsUser.Cells.Item(2, 3).Value = "ID tài khoản"

# Renumber the "STT" column for all the rows that got pushed down.
for ($r = 3; $r -le 14; $r++) {
    $wsUser.Cells.Item($r, 1).Value = $r - 1
}

# ---------------------------------------------------------------------------
# 2. Rename the worksheets, dropping the "Table_" prefix.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Table_Book").Name = "Book"
$wb.Worksheets.Item("Table_Author").Name = "Author"
$wb.Worksheets.Item("Table_Category").Name = "Category"
$wb.Worksheets.Item("Table_User").Name = "User"
$wb.Worksheets.Item("Table_Comment").Name = "Comment"

# ---------------------------------------------------------------------------
# 3. Update the selections / active sheet: the "User" sheet is no longer the
#    active tab (selection parked at C6); "Comment" becomes the active tab
#    (selection at C12, the spot the User tab used to have selected).
# ---------------------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("User")
$wsUser.Activate()
$wsUser.Range("C6").Select()

$wsComment = $wb.Worksheets.Item("Comment")
$wsComment.Activate()
$wsComment.Range("C12").Select()
